$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to column U ("Flow- Summer Base Flow") values per commit diff
$ws.Range("U3").Value = 1
$ws.Range("U10").Value = 3
$ws.Range("U12").Value = 1
$ws.Range("U13").Value = 1
$ws.Range("U14").Value = 1
$ws.Range("U15").Value = 1
$ws.Range("U16").Value = 1
$ws.Range("U17").Value = 1
$ws.Range("U18").Value = 1
$ws.Range("U19").Value = 1
$ws.Range("U20").Value = 1
$ws.Range("U23").Value = 1
$ws.Range("U28").Value = 1
$ws.Range("U30").Value = 3
$ws.Range("U31").Value = 3
$ws.Range("U32").Value = 3
$ws.Range("U33").Value = 3
$ws.Range("U34").Value = 3
$ws.Range("U35").Value = 3
$ws.Range("U41").Value = 3
$ws.Range("U45").Value = 3
$ws.Range("U47").Value = 3
$ws.Range("U52").Value = 3
$ws.Range("U55").Value = 5
$ws.Range("U56").Value = 3
$ws.Range("U57").Value = 5
$ws.Range("U58").Value = 5
$ws.Range("U59").Value = 5
$ws.Range("U60").Value = 5
$ws.Range("U61").Value = 5
$ws.Range("U62").Value = 5
$ws.Range("U63").Value = 3
$ws.Range("U64").Value = 5
$ws.Range("U65").Value = 5
$ws.Range("U66").Value = 5
$ws.Range("U68").Value = 5
$ws.Range("U69").Value = 5
$ws.Range("U70").Value = 1
$ws.Range("U71").Value = 5
$ws.Range("U72").Value = 3
$ws.Range("U73").Value = 3
$ws.Range("U74").Value = 3
$ws.Range("U75").Value = 5
$ws.Range("U76").Value = 3
$ws.Range("U79").Value = 1
$ws.Range("U80").Value = 3
$ws.Range("U81").Value = 1
$ws.Range("U82").Value = 3
$ws.Range("U83").Value = 3
$ws.Range("U84").Value = 3
$ws.Range("U85").Value = 1
$ws.Range("U86").Value = 1
$ws.Range("U88").Value = 1
$ws.Range("U89").Value = 1
$ws.Range("U90").Value = 1
$ws.Range("U91").Value = 1
$ws.Range("U92").Value = 1
$ws.Range("U93").Value = 1
$ws.Range("U94").Value = 1
$ws.Range("U115").Value = 1
$ws.Range("U116").Value = 3
$ws.Range("U117").Value = 1
$ws.Range("U119").Value = 3
$ws.Range("U120").Value = 3
$ws.Range("U121").Value = 1
$ws.Range("U122").Value = 1
$ws.Range("U123").Value = 1
$ws.Range("U124").Value = 1
$ws.Range("U127").Value = 1
$ws.Range("U128").Value = 1
$ws.Range("U129").Value = 1
